# Applies the cryptos-list refresh described by the commit diff:
# - Price (column D) and Volume/1h (column E) updates for existing rows
# - Rows 49-50 content swap (ONDO <-> Monero, re-ranked)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings (e.g. "6.30", "67.766.85"). Assigning a
# plain numeric-looking string via .Value lets Excel auto-convert it to a
# real number (dropping formatting / trailing zeros and adding a number
# style). Forcing the cell to Text format first, writing the literal
# string, then resetting the style back to "Normal" keeps the text exactly
# as authored while leaving the cell style untouched (same as original).
function Set-TextValue($cellRef, $val) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue 'D2' '67.766.85'
$ws.Range('E2').Value = '  +0.54%  '
Set-TextValue 'D3' '3.801.64'
$ws.Range('E3').Value = '  +0.48%  '
$ws.Range('E4').Value = '  +0.10%  '
Set-TextValue 'D5' '596.78'
$ws.Range('E5').Value = '  +0.56%  '
Set-TextValue 'D6' '167.29'
$ws.Range('E6').Value = '  +0.90%  '
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('E8').Value = '  +0.43%  '
$ws.Range('E9').Value = '  +1.58%  '
Set-TextValue 'D10' '6.30'
$ws.Range('E10').Value = '  -1.15%  '
$ws.Range('E11').Value = '  +0.30%  '
$ws.Range('E12').Value = '  -0.18%  '
Set-TextValue 'D13' '35.89'
$ws.Range('E13').Value = '  +0.32%  '
Set-TextValue 'D14' '4.443.28'
$ws.Range('E14').Value = '  +0.55%  '
Set-TextValue 'D15' '3.825.33'
$ws.Range('E15').Value = '  +1.15%  '
Set-TextValue 'D16' '18.55'
$ws.Range('E16').Value = '  +3.28%  '
Set-TextValue 'D17' '67.806.31'
$ws.Range('E17').Value = '  +0.66%  '
Set-TextValue 'D18' '7.07'
$ws.Range('E18').Value = '  +1.77%  '
$ws.Range('E19').Value = '  +0.64%  '
Set-TextValue 'D20' '461.13'
$ws.Range('E20').Value = '  +0.78%  '
Set-TextValue 'D21' '9.90'
$ws.Range('E21').Value = '  -2.88%  '
$ws.Range('E22').Value = '  -0.19%  '
$ws.Range('E23').Value = '  +1.22%  '
Set-TextValue 'D24' '83.29'
$ws.Range('E24').Value = '  -0.10%  '
$ws.Range('E25').Value = '  +2.02%  '
$ws.Range('E26').Value = '  -1.29%  '
$ws.Range('E27').Value = '  -0.05%  '
Set-TextValue 'D28' '9.99'
$ws.Range('E28').Value = '  +0.47%  '
Set-TextValue 'D29' '3.946.91'
$ws.Range('E29').Value = '  +0.40%  '
$ws.Range('E30').Value = '  -0.18%  '
Set-TextValue 'D31' '7.36'
$ws.Range('E31').Value = '  +2.50%  '
$ws.Range('E32').Value = '  +1.98%  '
Set-TextValue 'D33' '29.54'
$ws.Range('E33').Value = '  -0.77%  '
Set-TextValue 'D34' '1.00'
$ws.Range('E34').Value = '  +0.00%  '
Set-TextValue 'D35' '9.05'
$ws.Range('E35').Value = '  -1.06%  '
Set-TextValue 'D36' '3.741.75'
$ws.Range('E36').Value = '  +0.10%  '
$ws.Range('E37').Value = '  +0.54%  '
Set-TextValue 'D38' '3.36'
$ws.Range('E38').Value = '  +2.77%  '
$ws.Range('E39').Value = '  +0.13%  '
$ws.Range('E40').Value = '  +0.69%  '
$ws.Range('E41').Value = '  +1.13%  '
$ws.Range('E42').Value = '  +0.04%  '
$ws.Range('E43').Value = '  -0.01%  '
$ws.Range('E44').Value = '  +2.48%  '
Set-TextValue 'D45' '0.301'
$ws.Range('E45').Value = '  +1.27%  '
Set-TextValue 'D46' '42.73'
$ws.Range('E46').Value = '  -2.62%  '
Set-TextValue 'D47' '8.34'
$ws.Range('E47').Value = '  +0.04%  '
Set-TextValue 'D48' '27.27'
$ws.Range('E48').Value = '  +7.31%  '
$ws.Range('B49').Value = 'Monero'
$ws.Range('C49').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue 'D49' '147.49'
$ws.Range('E49').Value = '  -0.02%  '
$ws.Range('B50').Value = 'ONDO'
$ws.Range('C50').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
Set-TextValue 'D50' '1.36'
$ws.Range('E50').Value = '  +9.46%  '
Set-TextValue 'D51' '395.01'
$ws.Range('E51').Value = '  +0.60%  '
